$wb = $excel.ActiveWorkbook

# --- Sheet "safety_orders": remove the first safety order row (old row 2, Safety Order No. = 3) ---
$wsSafety = $wb.Worksheets.Item("safety_orders")
$wsSafety.Rows.Item(2).Delete()

# --- Sheet "open_buy_orders": add a new row for an additional buy order ---
$wsBuy = $wb.Worksheets.Item("open_buy_orders")
$wsBuy.Range("A4").Value = "OMNIAC-TCBSW-DM6PM5"
$wsBuy.Range("B4").Value = 158.09
